# Auto-generated Excel COM-interop script
# Applies updated 'F' (want-to-go count) and 'G' (lowest ticket price) values
# as described by the commit diff (bot-refreshed scrape data).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("G4").Value = 90
$ws.Range("F5").Value = 1087
$ws.Range("F6").Value = 28
$ws.Range("F7").Value = 133
$ws.Range("F8").Value = 1389
$ws.Range("F9").Value = 59
$ws.Range("F10").Value = 85
$ws.Range("F11").Value = 400
$ws.Range("F12").Value = 120
$ws.Range("F13").Value = 74
$ws.Range("F15").Value = 433
$ws.Range("F16").Value = 459
$ws.Range("F17").Value = 124
$ws.Range("F18").Value = 23
$ws.Range("F19").Value = 529
$ws.Range("F20").Value = 2526
$ws.Range("F22").Value = 37
$ws.Range("F28").Value = 92
$ws.Range("F30").Value = 910
$ws.Range("F32").Value = 27
$ws.Range("F36").Value = 231
$ws.Range("F37").Value = 21
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 709
$ws.Range("F5").Value = 589
$ws.Range("F6").Value = 589
$ws.Range("F15").Value = 290
$ws.Range("F16").Value = 290
$ws.Range("F19").Value = 920
$ws.Range("F22").Value = 594
$ws.Range("F26").Value = 211
$ws.Range("F27").Value = 216
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 71
$ws.Range("F6").Value = 2250
$ws.Range("F7").Value = 887
$ws.Range("F10").Value = 1081
$ws.Range("F12").Value = 67
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 71
$ws.Range("F6").Value = 2250
$ws.Range("F11").Value = 887
$ws.Range("F12").Value = 1081
$ws.Range("G13").Value = 90
$ws.Range("F15").Value = 67
$ws.Range("F16").Value = 709
$ws.Range("F17").Value = 1087
$ws.Range("F18").Value = 133
$ws.Range("F19").Value = 589
$ws.Range("F20").Value = 59
$ws.Range("F21").Value = 85
$ws.Range("F22").Value = 400
$ws.Range("F23").Value = 120
$ws.Range("F24").Value = 74
$ws.Range("F26").Value = 433
$ws.Range("F27").Value = 459
$ws.Range("F28").Value = 124
$ws.Range("F29").Value = 23
$ws.Range("F30").Value = 529
$ws.Range("F36").Value = 92
$ws.Range("F39").Value = 910
$ws.Range("F40").Value = 290
$ws.Range("F44").Value = 211
$ws.Range("F49").Value = 231
$ws.Range("F50").Value = 21
